$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 20
$ws.Range("A5").Value = 50

$ws.Range("A5").Select()
